$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the polyA isolation kit (NEB kit code gained an "L" suffix)
$ws.Range("G2:G27").Value = "NEBNextPoly(A)E7490L"

# roboticS1Prep column: replace the "No" text answers with a real boolean
# FALSE value, formatted to display as TRUE/FALSE
$ws.Range("I2:I27").Value = $false
$ws.Range("I2:I27").NumberFormat = '"TRUE";"TRUE";"FALSE"'

# Widen the polyA isolation protocol column so the longer kit name fits
$ws.Columns("G").ColumnWidth = 34.6

# Move the active selection to the column that was just edited
$ws.Range("G2:G27").Select() | Out-Null
